# Weekly update: a new price observation for "Choclo" (Comercializadora del
# Agro de Limarí) is inserted as a new record. The new record becomes the
# first data row (row 68), pushing the existing data rows down by one
# (old row 68 -> new row 69, ..., old row 154 -> new row 155).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 68; everything currently at/after row 68
# (data rows through the old last row 154) shifts down by one row.
$ws.Rows.Item(68).Insert()

# Populate the newly inserted row 68 with the new record's data.
$ws.Range("A68").Value = 2
$ws.Range("B68").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C68").Value = "Coquimbo"
$ws.Range("D68").Value = 44944
$ws.Range("E68").Value = 4
$ws.Range("F68").Value = 100112024
$ws.Range("G68").Value = "Choclo"
$ws.Range("H68").Value = "Choclero"
$ws.Range("I68").Value = "Primera"
$ws.Range("J68").Value = 110000
$ws.Range("K68").Value = 200
$ws.Range("L68").Value = 250
$ws.Range("M68").Value = 225
$ws.Range("N68").Value = "`$/unidad"
$ws.Range("O68").Value = "Provincia de Limarí"
$ws.Range("P68").Value = 225
$ws.Range("Q68").Value = 1
$ws.Range("R68").Value = "Hortaliza"
